# "played with weights, updated spreadsheet"
#
# Adjust several Numeric Rating weights, add a note explaining the -1
# "loves Jesus" sentinel value, rename one Importance entry, and insert
# three new criteria rows ("wants kids", a blank spacer, "genuine")
# just above the closing "works out" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric / text value updates (rows 2-35; unaffected by the later row insert) ---
$ws.Range("C2").Value = 35
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 25
$ws.Range("C6").Value = 25
$ws.Range("C7").Value = 1
$ws.Range("C11").Value = 5
$ws.Range("D12").Value = "the -1 for parsing"
$ws.Range("C14").Value = 20
$ws.Range("C16").Value = 45
$ws.Range("C17").Value = 40
$ws.Range("C22").Value = 20
$ws.Range("C23").Value = 30
$ws.Range("C25").Value = 25
$ws.Range("B33").Value = "encouraged"

# --- Insert three new rows above the former last data row (row 36, "works out") ---
[void]$ws.Rows.Item(36).Insert()
[void]$ws.Rows.Item(36).Insert()
[void]$ws.Rows.Item(36).Insert()

# Copy formatting from row 35 ("bright eyes") down onto the three fresh rows so they
# inherit the same cell styles used by the rest of the table.
[void]$ws.Range("A35:D35").Copy()
[void]$ws.Range("A36:D38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 36: wants kids
$ws.Range("A36").Value = "wants kids"
$ws.Range("B36").Value = "important"
$ws.Range("C36").Value = 25
$ws.Range("D36").Value = "Psalm 127:4"

# Row 37: left blank (spacer row)

# Row 38: genuine
$ws.Range("A38").Value = "genuine"
$ws.Range("B38").Value = "essential"
$ws.Range("C38").Value = 40

# --- Conditional formatting on the "Red flags" notes column shifted down with the
# --- inserted rows (old D44 -> new D47); re-create it at the new location.
$ws.Range("D44").FormatConditions.Delete()
$cfRng = $ws.Range("D47")
$cf = $cfRng.FormatConditions.Add(9, 3, "high")
$cf.Formula1 = "=NOT(ISERROR(SEARCH(""high"",D47)))"

# --- View cosmetics: zoom + selection ---
$excel.ActiveWindow.Zoom = 150
[void]$ws.Range("A37").Select()

Write-Output "edit complete"
